$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 11: a work-log entry (date + quantity of hours), continuing the
# existing table in columns A:B.
$ws.Range("A11").Value = 41563
$ws.Range("B11").Value = 0.09375

# Reuse the same number formats already used by the rest of the table
# (A: date; B: time) instead of letting Excel mint brand-new style
# records for the new cells.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# Move the active selection to the cell right after the new row, mirroring
# how the sheet was left selected after the edit.
$ws.Range("C11").Select()
